{"js": "// The <id> paragraph currently holds three runs:\n//   \"<id>\"  (Courier New, color 7f6000, sz 18)\n//   \"p030r_a2\"  (plain formatting)\n//   \"</id>\" (Courier New, color 7f6000, sz 18)\n// The edit collapses them into a single run reading \"<id>p030r_2</id>\"\n// (using the formatting of the first run), i.e. the old identifier\n// \"p030r_a2\" is replaced by \"p030r_2\" and the surrounding tag runs are\n// merged with it into one run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText = \"<id>p030r_a2</id>\";\nconst newText = \"<id>p030r_2</id>\";\n\nlet target = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === oldText) {\n    target = paragraph;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the '<id>p030r_a2</id>' paragraph.\");\n}\n\n// Replacing the whole paragraph range's text in one shot merges the\n// paragraph's runs into a single run that carries the formatting of the\n// paragraph's first run \u2014 matching the target OOXML (one <w:r> with the\n// Courier New / 7f6000 / sz18 rPr).\ntarget.getRange().insertText(newText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# The \"<id>\" paragraph currently holds three runs:\n#   \"<id>\"      (Courier New, color 7f6000, sz 18)\n#   \"p030r_a2\"  (plain formatting)\n#   \"</id>\"     (Courier New, color 7f6000, sz 18)\n# The edit collapses them into a single run reading \"<id>p030r_2</id>\"\n# (carrying the formatting of the first run) - i.e. the identifier\n# \"p030r_a2\" becomes \"p030r_2\" and the three runs merge into one.\n\n$d = $word.ActiveDocument\n\n$oldText = \"<id>p030r_a2</id>\"\n$newText = \"<id>p030r_2</id>\"\n\n# Find.Execute repositions/collapses the range onto the matched text\n# (paragraph mark excluded), so setting .Text on it afterwards replaces\n# exactly \"<id>p030r_a2</id>\" and lets Word re-flow the paragraph's runs\n# into one, inheriting the first run's character formatting.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($oldText, $false, $true)\n\nif ($found) {\n    $searchRange.Text = $newText\n} else {\n    # Fallback: locate the paragraph by its exact text and replace its\n    # content (excluding the trailing paragraph mark) directly.\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $oldText) {\n            $target = $p\n            break\n        }\n    }\n    if ($target -eq $null) {\n        throw \"Could not find the '<id>p030r_a2</id>' paragraph.\"\n    }\n    $r = $target.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $newText\n}\n"}
